$d = $word.ActiveDocument

# 1. Switch the licence text from CC BY-NC to CC BY-SA in the footer
#    paragraph ("Except where otherwise noted ... licensed under CC BY-NC 4.0 ...").
$d.Content.Find.Execute("CC BY-NC 4.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CC BY-SA 4.0", 2) | Out-Null

# 2. Point the "https://creativecommons.org/licenses/by-nc/4.0" hyperlink (both
#    its visible text and its target address) at the by-sa licence instead.
$h = $d.Hyperlinks(1)
$h.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
$h.Address = "https://creativecommons.org/licenses/by-sa/4.0"

# 3. Re-touch the "Follow the instructor directions" bullet so the stale
#    lastRenderedPageBreak hint Word had cached for it is dropped.
$d.Content.Find.Execute("Follow the instructor directions", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "Follow the instructor directions", 2) | Out-Null
